# Add CBO e faz ajustes nas pages de MI e MC
#
# 1. Rename the existing "Concepts" sheet to "Properties" and replace its
#    content with the CodeSystem property table (Code/Uri/Description/Type).
# 2. Add a brand-new "Concepts" sheet (after "Properties") with the original
#    concept rows, plus:
#      - a new "0" / "Sem registro no sistema de informação de origem" row
#      - updated display text for code "9"
#      - a new "15" / "Operação Gota" row

$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# ---------------------------------------------------------------------
# Step 1: turn the old "Concepts" sheet into the new "Properties" sheet
# ---------------------------------------------------------------------
$props = $wb.Worksheets.Item("Concepts")
$props.Name = "Properties"

# Drop the old concept rows (3-15); only 2 rows are needed going forward.
$props.Range("A3:D15").Clear()

$props.Range("A1").Value = "Code"
$props.Range("B1").Value = "Uri"
$props.Range("C1").Value = "Description"
$props.Range("D1").Value = "Type"

$props.Range("A2").Value = "inactive"
$props.Range("B2").Value = "http://hl7.org/fhir/concept-properties#inactive"
$props.Range("C2").Value = ""
$props.Range("D2").Value = "boolean"

# ---------------------------------------------------------------------
# Step 2: create the new "Concepts" sheet after "Properties"
# ---------------------------------------------------------------------
$concepts = $wb.Worksheets.Add($null, $props)
$concepts.Name = "Concepts"

$concepts.Range("A1").Value = "Level"
$concepts.Range("B1").Value = "Code"
$concepts.Range("C1").Value = "Display"
$concepts.Range("D1").Value = "Definition"

# Match the look & feel of the other sheets: header style from row 1,
# bordered body style from row 2 onward.
$meta.Range("A1:B1").Copy()
$concepts.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Range("A2:B2").Copy()
$concepts.Range("A2:D17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @("0", "Sem registro no sistema de informação de origem"),
    @("1", "Rotina"),
    @("2", "Especial"),
    @("3", "Bloqueio"),
    @("4", "Intensificação"),
    @("5", "Campanha indiscriminada"),
    @("6", "Campanha seletiva"),
    @("7", "Soroterapia"),
    @("8", "Serviço Privado"),
    @("9", "Monitoramento das Estratégias de Vacinação"),
    @("10", "Pesquisa"),
    @("11", "Pré-exposição"),
    @("12", "Pós-exposição"),
    @("13", "Reexposição"),
    @("14", "Vacinação Escolar"),
    @("15", "Operação Gota")
)

# Force the "Level"/"Code" columns to store text (e.g. "0", "1", ...)
# rather than being auto-converted to numbers, matching the source data.
$concepts.Range("A2").NumberFormat = "@"
$concepts.Range("A2").Value = "1"
$concepts.Range("B2").NumberFormat = "@"
$concepts.Range("B2").Value = "0"
$concepts.Range("A2:B2").Copy()
$concepts.Range("A2:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$r = 2
foreach ($row in $rows) {
    $concepts.Range("A$r").Value = "1"
    $concepts.Range("B$r").Value = $row[0]
    $concepts.Range("C$r").Value = $row[1]
    $r = $r + 1
}
